$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two new data rows at row 569, pushing the existing rows 569-593
# down to 571-595 (weekly update: newest week's prices prepended to the
# top of this market/product block).
$ws.Rows.Item(569).Insert()
$ws.Rows.Item(569).Insert()

# New row 569: Primera quality, 2023-05-29
$ws.Range("A569").Value = 4
$ws.Range("B569").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C569").Value = "Los Lagos"
$ws.Range("D569").Value = "5/29/2023"
$ws.Range("E569").Value = 10
$ws.Range("F569").Value = "Fruta"
$ws.Range("G569").Value = 100102
$ws.Range("H569").Value = "Cítricos"
$ws.Range("I569").Value = 100102006
$ws.Range("J569").Value = "Pomelo"
$ws.Range("K569").Value = "Start Ruby"
$ws.Range("L569").Value = "Primera"
$ws.Range("M569").Value = 80
$ws.Range("N569").Value = 14000
$ws.Range("O569").Value = 15000
$ws.Range("P569").Value = 14500
$ws.Range("Q569").Value = "$/caja 14 kilos empedrada"
$ws.Range("R569").Value = "Región de O'Higgins"
$ws.Range("S569").Value = 1036
$ws.Range("T569").Value = 14

# New row 570: Segunda quality, 2023-05-29
$ws.Range("A570").Value = 4
$ws.Range("B570").Value = "Feria Lagunitas de Puerto Montt"
$ws.Range("C570").Value = "Los Lagos"
$ws.Range("D570").Value = "5/29/2023"
$ws.Range("E570").Value = 10
$ws.Range("F570").Value = "Fruta"
$ws.Range("G570").Value = 100102
$ws.Range("H570").Value = "Cítricos"
$ws.Range("I570").Value = 100102006
$ws.Range("J570").Value = "Pomelo"
$ws.Range("K570").Value = "Start Ruby"
$ws.Range("L570").Value = "Segunda"
$ws.Range("M570").Value = 40
$ws.Range("N570").Value = 13000
$ws.Range("O570").Value = 13000
$ws.Range("P570").Value = 13000
$ws.Range("Q570").Value = "$/caja 14 kilos empedrada"
$ws.Range("R570").Value = "Región de O'Higgins"
$ws.Range("S570").Value = 929
$ws.Range("T570").Value = 14
